$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 211, pushing existing rows 211+ down to 213+
$ws.Rows("211:212").Insert()

# Fill in row 211 (new weekly data point: Zafiro rojo)
$ws.Range("A211").Value = 7
$ws.Range("B211").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C211").Value = "Ñuble"
$ws.Range("D211").Value = 44795
$ws.Range("E211").Value = 16
$ws.Range("F211").Value = 100112002
$ws.Range("G211").Value = "Pimiento"
$ws.Range("H211").Value = "Zafiro rojo"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 120
$ws.Range("K211").Value = 26000
$ws.Range("L211").Value = 27000
$ws.Range("M211").Value = 26500
$ws.Range("N211").Value = "$/caja 15 kilos"
$ws.Range("O211").Value = "Región de Arica y Parinacota"
$ws.Range("P211").Value = 1767
$ws.Range("Q211").Value = 15
$ws.Range("R211").Value = "Hortaliza"

# Fill in row 212 (new weekly data point: Zafiro verde)
$ws.Range("A212").Value = 7
$ws.Range("B212").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C212").Value = "Ñuble"
$ws.Range("D212").Value = 44795
$ws.Range("E212").Value = 16
$ws.Range("F212").Value = 100112002
$ws.Range("G212").Value = "Pimiento"
$ws.Range("H212").Value = "Zafiro verde"
$ws.Range("I212").Value = "Primera"
$ws.Range("J212").Value = 80
$ws.Range("K212").Value = 25000
$ws.Range("L212").Value = 25000
$ws.Range("M212").Value = 25000
$ws.Range("N212").Value = "$/caja 15 kilos"
$ws.Range("O212").Value = "Región de Arica y Parinacota"
$ws.Range("P212").Value = 1667
$ws.Range("Q212").Value = 15
$ws.Range("R212").Value = "Hortaliza"
